{"js": "// Load the paragraphs so we can find the two anchor paragraphs by their text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet lukeParagraph = null;\nlet assuntosParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = p.text.trim();\n  if (text.indexOf(\"Luke Marques Magalh\u00e3es\") !== -1 && lukeParagraph === null) {\n    lukeParagraph = p;\n  }\n  if (text === \"Assuntos abordados:\" && assuntosParagraph === null) {\n    assuntosParagraph = p;\n  }\n}\n\n// 1) New paragraph right after \"Luke Marques Magalh\u00e3es \u2013 32218605\" with the\n//    second team member's name, following the \"Normal\" style.\nif (lukeParagraph) {\n  const novoIntegrante = lukeParagraph.insertParagraph(\n    \"Luis Augusto Marques - 32237820\",\n    \"After\"\n  );\n  novoIntegrante.style = \"Normal\";\n}\n\n// 2) New bold paragraph right after \"Assuntos abordados:\" introducing the topic.\nif (assuntosParagraph) {\n  const intro = assuntosParagraph.insertParagraph(\n    \"Em nosso projeto, utilizamos diversos conte\u00fados abordados em sala, neste arquivo iremos comentar sobre o uso dos mesmos.\",\n    \"After\"\n  );\n  intro.font.bold = true;\n}\n\n// 3) Make the page orientation explicitly portrait (already portrait-sized,\n//    but the section properties should carry the explicit attribute).\nconst section = context.document.sections.getFirst();\nsection.pageSetup.orientation = Word.PageOrientation.portrait;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- helper: find a paragraph whose text matches a -like pattern ---\nfunction Find-Paragraph($doc, $pattern) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text -like $pattern) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# 1) Insert a new paragraph right after \"Luke Marques Magalhaes - 32218605\"\n#    with the second team member's name, using the \"Normal\" style.\n$lukeParagraph = Find-Paragraph $d \"*Luke Marques*32218605*\"\nif ($lukeParagraph -ne $null) {\n    $rng = $lukeParagraph.Range\n    $rng.Collapse(0)          # wdCollapseEnd\n    $rng.InsertParagraphAfter()\n    $novoIntegrante = $lukeParagraph.Next()\n    $novoIntegrante.Range.Text = \"Luis Augusto Marques - 32237820\"\n    $novoIntegrante.Range.Style = \"Normal\"\n}\n\n# 2) Insert a new bold paragraph right after \"Assuntos abordados:\" introducing\n#    the topic. Re-locate the paragraph fresh (the document just mutated, so\n#    any previously captured paragraph reference may be stale).\n$assuntosParagraph = Find-Paragraph $d \"Assuntos abordados:*\"\nif ($assuntosParagraph -ne $null) {\n    $rng2 = $assuntosParagraph.Range\n    $rng2.Collapse(0)         # wdCollapseEnd\n    $rng2.InsertParagraphAfter()\n    $intro = $assuntosParagraph.Next()\n    $intro.Range.Text = \"Em nosso projeto, utilizamos diversos conte\u00fados abordados em sala, neste arquivo iremos comentar sobre o uso dos mesmos.\"\n    $intro.Range.Font.Bold = 1\n}\n\n# 3) Make the page orientation explicitly portrait.\n$d.PageSetup.Orientation = 0   # wdOrientPortrait\n"}
